# Add new dataset rows (22Q2 EH numbers) to depmap_datasets_list sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rsquo = [char]0x2019

# New rows to append starting at row 78
$rows = @(
    @("crispr_22Q2", "(Chronos) Batch and off-target corrected CRISPR-Cas9 gene knockout dependency data", "EH7554", 17386, 1086, 31, 28, 3.16, "May 9 2022", "https://ndownloader.figshare.com/files/34990036"),
    @("copyNumber_22Q2", "Inferred log copy number data", "EH7555", 25368, 1766, 33, 30, 3.16, "May 9 2022", "https://ndownloader.figshare.com/files/34989937"),
    @("TPM_22Q2", "CCLE TPM RNAseq gene expression data for protein coding genes", "EH7556", 19221, 1406, 33, 30, 3.16, "May 9 2022", "https://ndownloader.figshare.com/files/34989919"),
    @("mutationCalls_22Q2", "Merged mutation calls (for coding region, germline filtered) ", "EH7557", 18784, 1771, 33, 30, 3.16, "May 9 2022", "https://ndownloader.figshare.com/files/34989940"),
    @("metadata_22Q2", "Metadata for cell lines in the 22Q1 DepMap release", "EH7558", "NA", 1840, 33, 30, 3.16, "May 9 2022", "https://ndownloader.figshare.com/files/35020903"),
    @("achilles_22Q2", "Project Achilles$rsquo CRISPR screen metadata", "EH7559", "NA", "NA", "NA", "NA", 3.16, "May 9 2022", "https://ndownloader.figshare.com/files/34989901")
)

$startRow = 78
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 9).NumberFormat = "@"
    $ws.Cells.Item($r, 10).Value = $row[9]
}

# Update selection to mirror diff (topLeftCell A63, activeCell B87)
$ws.Application.ActiveWindow.ScrollRow = 63
$ws.Range("B87").Select()
